$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before AJ (old AJ/AK shift right to AK/AL)
$ws.Range("AJ1").EntireColumn.Insert()

# New header + value for the inserted column (now AJ)
$ws.Range("AJ1").Value = "MgCa Coretop modelled temperature"
$ws.Range("AJ2").Value = 13.8771

# Update the values that changed in row 2 for existing columns
$ws.Range("X2").Value = 13.35
$ws.Range("Y2").Value = 4.399195263974608
$ws.Range("Z2").Value = 3.70739198397461
$ws.Range("AA2").Value = -1.68524917602539
$ws.Range("AB2").Value = -0.8866664487526901

# The previously-shifted AL2 value has a tiny precision tweak
$ws.Range("AL2").Value = -1.418517270000001
